$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 02:07"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 7287421
$ws.Range("C4").Value = 43066
$ws.Range("D4").Value = 4519379
$ws.Range("E4").Value = 2558865
$ws.Range("G4").Value = 737
$ws.Range("H4").Value = 209177

# Brasil (row 6)
$ws.Range("B6").Value = 4718115
$ws.Range("C6").Value = 25536
$ws.Range("E6").Value = 525837
$ws.Range("G6").Value = 732
$ws.Range("H6").Value = 141441

# Barein (row 56)
$ws.Range("B56").Value = 68775
$ws.Range("C56").Value = 585
$ws.Range("D56").Value = 62252
$ws.Range("E56").Value = 6284

# Chequia (row 57)
$ws.Range("B57").Value = 63294
$ws.Range("C57").Value = 1976
$ws.Range("D57").Value = 30936
$ws.Range("E57").Value = 31767
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 591

# Rows 71-72 swap order: Paraguay now ranks above Kenia, with refreshed data
$ws.Range("A71").Value = "Paraguay"
$ws.Range("B71").Value = 37922
$ws.Range("C71").Value = 696
$ws.Range("D71").Value = 21757
$ws.Range("E71").Value = 15383
$ws.Range("G71").Value = 21
$ws.Range("H71").Value = 782

$ws.Range("A72").Value = "Kenia"
$ws.Range("B72").Value = 37871
$ws.Range("C72").Value = 164
$ws.Range("D72").Value = 24581
$ws.Range("E72").Value = 12601
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 689

# Camerun (row 83)
$ws.Range("B83").Value = 20735
$ws.Range("C83").Value = 23
$ws.Range("E83").Value = 877

# Bulgaria (row 84)
$ws.Range("B84").Value = 19997
$ws.Range("C84").Value = 169
$ws.Range("D84").Value = 14160
$ws.Range("E84").Value = 5048

# Maldivas (row 101)
$ws.Range("B101").Value = 10045
$ws.Range("C101").Value = 31
$ws.Range("D101").Value = 8754
$ws.Range("E101").Value = 1257

# Mauritania (row 113)
$ws.Range("B113").Value = 7462
$ws.Range("C113").Value = 5
$ws.Range("E113").Value = 231

# Surinam (row 127)
$ws.Range("B127").Value = 4831
$ws.Range("C127").Value = 14
$ws.Range("D127").Value = 4620
$ws.Range("E127").Value = 109

# Uruguay (row 154)
$ws.Range("B154").Value = 1998
$ws.Range("C154").Value = 31
$ws.Range("D154").Value = 1716
$ws.Range("E154").Value = 235

# Rows 206-207 swap order: Santa Lucia now ranks above Timor Oriental (values identical, only labels swap)
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("A207").Value = "Timor Oriental"
